# App_pressing_Loreau / LogPattern.xlsx
# "Version Non fonctionnelle : génération excel"
#
# Add the header label used when the log count is written into the sheet,
# size the column to fit it, and leave the selection where the user last
# clicked (D6) before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "Nombre de logs"
$ws.Columns.Item(1).ColumnWidth = 16.140625

[void]$ws.Range("D6").Select()
